# add staff id and phone number to line manager
#
# Makes room for a new "line manager" column by shifting the existing
# secondary-phone-number column (G) one column to the right (H) on every
# data row. G is left blank (but keeps its formatting) ready to receive
# the new staff id / phone number data for the line manager.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move G1:G10 -> H1:H10 (values, formatting and shared-string identity all
# travel with the cut/paste, and the source cells are cleared automatically).
$ws.Range("G1:G10").Cut($ws.Range("H1:H10"))

# Keep the new column's width consistent with the rest of the table so the
# worksheet's column metadata covers the newly used column H.
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
